$wb = $excel.ActiveWorkbook

# ----- Summary sheet -----
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1200.77
$wsSummary.Range("B4").Value = 0.78
$wsSummary.Range("B6").Value = 122
$wsSummary.Range("B8").Value = 46
$wsSummary.Range("B9").Value = 45.08

# ----- Strategy Status sheet -----
$wsStrategy = $wb.Worksheets.Item("Strategy Status")
$wsStrategy.Range("C4").Value = 100.77
$wsStrategy.Range("D4").Value = 122
$wsStrategy.Range("E4").Value = 0.78
$wsStrategy.Range("F4").Value = 0.77
$wsStrategy.Range("G4").Value = 45.08

# ----- All Trades sheet: append new closed trade (#122) -----
$wsAllTrades = $wb.Worksheets.Item("All Trades")
$wsAllTrades.Range("A123").Value = 122
$wsAllTrades.Range("B123").NumberFormat = "@"
$wsAllTrades.Range("B123").Value = "2026-02-17"
$wsAllTrades.Range("B123").Style = "Normal"
$wsAllTrades.Range("C123").NumberFormat = "@"
$wsAllTrades.Range("C123").Value = "09:29:09"
$wsAllTrades.Range("C123").Style = "Normal"
$wsAllTrades.Range("D123").Value = "MarketMaking"
$wsAllTrades.Range("E123").Value = "DOWN"
$wsAllTrades.Range("F123").Value = 0.9399999999999999
$wsAllTrades.Range("G123").Value = 0.92
$wsAllTrades.Range("H123").Value = "CLOSED"
$wsAllTrades.Range("I123").Value = -2.1277
$wsAllTrades.Range("J123").Value = -0.02
$wsAllTrades.Range("K123").Value = 100.77
$wsAllTrades.Range("L123").Value = 0
$wsAllTrades.Range("M123").Value = 0
$wsAllTrades.Range("N123").Value = 0.6
$wsAllTrades.Range("O123").Value = "Normal spread capture: 19600 bps"
$wsAllTrades.Range("P123").Value = "early_exit"
$wsAllTrades.Range("Q123").Value = 0.14

# ----- MarketMaking sheet: append the same closed trade (#122) -----
$wsMM = $wb.Worksheets.Item("MarketMaking")
$wsMM.Range("A123").Value = 122
$wsMM.Range("B123").NumberFormat = "@"
$wsMM.Range("B123").Value = "2026-02-17"
$wsMM.Range("B123").Style = "Normal"
$wsMM.Range("C123").NumberFormat = "@"
$wsMM.Range("C123").Value = "09:29:09"
$wsMM.Range("C123").Style = "Normal"
$wsMM.Range("D123").Value = "MarketMaking"
$wsMM.Range("E123").Value = "DOWN"
$wsMM.Range("F123").Value = 0.9399999999999999
$wsMM.Range("G123").Value = 0.92
$wsMM.Range("H123").Value = "CLOSED"
$wsMM.Range("I123").Value = -2.1277
$wsMM.Range("J123").Value = -0.02
$wsMM.Range("K123").Value = 100.77
$wsMM.Range("L123").Value = 0
$wsMM.Range("M123").Value = 0
$wsMM.Range("N123").Value = 0.6
$wsMM.Range("O123").Value = "Normal spread capture: 19600 bps"
$wsMM.Range("P123").Value = "early_exit"
$wsMM.Range("Q123").Value = 0.14
